# feat: add 2022-Q1 data
#
# 1) Duplicate the existing "总计" (grand-total) sheet, so the duplicate
#    keeps all of its original formatting (sheetPr / header styles /
#    page margins) to become the new grand-total sheet.
# 2) Turn the ORIGINAL "总计" sheet into the new "2022-Q1" per-fund detail
#    sheet (same column layout as the 2021-Q3 / 2021-Q4 sheets).
# 3) Turn the DUPLICATE sheet back into "总计" and prepend a 2022-Q1 row
#    to its summary table (date / holding count / holding value).

$wb = $excel.ActiveWorkbook

$q1 = $wb.Worksheets.Item("总计")
$q1.Copy($null, $q1)
$total = $wb.Worksheets.Item("总计 (2)")

# ---------------------------------------------------------------------
# Step 1: turn the old "总计" sheet into the new "2022-Q1" detail sheet
# ---------------------------------------------------------------------
$q1.Name = "2022-Q1"

# Extend the header formatting (currently only B1:D1 carry style) across
# the new columns E1:H1 so every header cell shares the same look.
$q1.Range("B1:D1").Copy()
$q1.Range("E1:H1").PasteSpecial(-4122)

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Row-index column (A) formatting: A2:A3 already carry style from the old
# sheet; copy it down onto the two new rows (A4:A5).
$q1.Range("A2:A3").Copy()
$q1.Range("A4:A5").PasteSpecial(-4122)

$q1.Range("A2").Value = 0
$q1.Range("A3").Value = 1
$q1.Range("A4").Value = 2
$q1.Range("A5").Value = 3

# Fund-code / name / ratio columns are stored as plain text in the source
# workbook (so leading zeros like "002567" and trailing zeros like
# "0.00" survive) -- pre-format the block as Text before writing. G5 is
# left out: it holds a genuine numeric 0, not formatted text.
$q1.Range("B2:G4").NumberFormat = "@"
$q1.Range("B5:F5").NumberFormat = "@"

$q1.Range("B2").Value = "002567"
$q1.Range("C2").Value = "大成国家安全主题灵活配置混合"
$q1.Range("D2").Value = "0.34"
$q1.Range("E2").Value = "52.90"
$q1.Range("F2").Value = "3.45"
$q1.Range("G2").Value = "0.0117"
$q1.Range("H2").Value = 9

$q1.Range("B3").Value = "002303"
$q1.Range("C3").Value = "金鹰智慧生活灵活配置混合"
$q1.Range("D3").Value = "0.11"
$q1.Range("E3").Value = "89.88"
$q1.Range("F3").Value = "3.06"
$q1.Range("G3").Value = "0.0034"
$q1.Range("H3").Value = 9

$q1.Range("B4").Value = "011444"
$q1.Range("C4").Value = "创金合信瑞裕混合A"
$q1.Range("D4").Value = "0.03"
$q1.Range("E4").Value = "68.73"
$q1.Range("F4").Value = "2.50"
$q1.Range("G4").Value = "0.0008"
$q1.Range("H4").Value = 8

$q1.Range("B5").Value = "011445"
$q1.Range("C5").Value = "创金合信瑞裕混合C"
$q1.Range("D5").Value = "0.00"
$q1.Range("E5").Value = "68.73"
$q1.Range("F5").Value = "2.50"
$q1.Range("G5").Value = 0
$q1.Range("H5").Value = 8

# ---------------------------------------------------------------------
# Step 2: turn the duplicate sheet back into "总计", with the new
# 2022-Q1 row inserted above the previously-existing rows.
# ---------------------------------------------------------------------
$total.Name = "总计"

$total.Rows(2).Insert()
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)
$total.Range("B2:D2").ClearFormats()

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 4
$total.Range("D2").Value = 0.02

# The index column keeps counting up after the insert (it is not the old
# values shifted down as-is).
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2

Write-Host "2022-Q1 sheet created and 总计 summary rebuilt"
